$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New applicant/officer inquiry row appended below the existing header row.
# Columns: A=ApplicantNRIC, B=Message, C=ProjectName, D=Reply, E=Timestamp
$ws.Range("A2").Value = "S1234567A"
$ws.Range("B2").Value = "fuck"
$ws.Range("C2").Value = "Acacia Breeze"

# Reply is still blank for this inquiry, but the cell itself must exist
# (touching a formatting property materializes it without changing style).
$ws.Range("D2").Font.Size = 11

# Timestamp stored as a date serial (20 Apr 2025), matching the existing
# dd/mm/yyyy date format already used in the workbook.
$ws.Range("E2").Value = 45767.0
$ws.Range("E2").NumberFormat = "dd/mm/yyyy"
